# Automatische test-sync: 2025-06-17 21:49:25
# Append two new mail-log rows to the "Logs" sheet and refresh the
# "Dashboard" summary counts accordingly.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Row 31: "Wat zijn jullie openingstijden?" (answered) ---
$logs.Range("A31").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B31").Value = "mailmind.test@zohomail.eu"
$logs.Range("C31").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D31").Value = "Informatieaanvraag"
$logs.Range("E31").Value = "Beste klant,`r`nBedankt voor je bericht. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 18:00 uur. Op zaterdag zijn we geopend van 10:00 tot 16:00 uur. Op zondag zijn we gesloten. Mocht je nog vragen hebben, dan hoor ik het graag.`r`nMet vriendelijke groet,`r`n[Naam] E-mailassistent"
$logs.Range("F31").Value = "2025-06-17 21:49:20"
$logs.Range("G31").Value = "Ja"
# Avoid a lingering auto-set custom row height from the multi-line E31 value.
$logs.Rows.Item(31).AutoFit()

# --- Row 32: "Afmelding nieuwsbrief" (not answered) ---
$logs.Range("A32").Value = "Afmelding nieuwsbrief"
$logs.Range("B32").Value = "mailmind.test@zohomail.eu"
$logs.Range("C32").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D32").Value = "Afmelding"
$logs.Range("F32").Value = "2025-06-17 21:49:21"
$logs.Range("G32").Value = "Nee"

# --- Extend conditional formatting ranges to cover the new rows ---
$logs.Range("D2:D30").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D32"))
$logs.Range("G2:G30").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G32"))

# --- Update Dashboard category counts ---
# Informatieaanvraag: 14 -> 15 (new row 31)
$dash.Range("B2").Value = 15
# Afmelding: 4 -> 5 (new row 32)
$dash.Range("B4").Value = 5
